$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 16 de Abril de 2020 a las 06:52"

# --- Row 23: India - update stats (no row swap) ---
$ws.Range("B23").Value = 12456
$ws.Range("C23").Value = 86
$ws.Range("D23").Value = 1513
$ws.Range("E23").Value = 10521

# --- Rows 34/35: Pakistan overtakes Australia, swap rows and update Pakistan's data ---
$ws.Range("A34").Value = "Pakistan"
$ws.Range("B34").Value = 6505
$ws.Range("C34").Value = 122
$ws.Range("D34").Value = 1645
$ws.Range("E34").Value = 4736
$ws.Range("F34").Value = 46
$ws.Range("G34").Value = 13
$ws.Range("H34").Value = 124

$ws.Range("A35").Value = "Australia"
$ws.Range("B35").Value = 6462
$ws.Range("C35").Value = 15
$ws.Range("D35").Value = 3702
$ws.Range("E35").Value = 2697
$ws.Range("F35").Value = 76
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 63

# --- Row 53: Tailandia - update stats (no row swap) ---
$ws.Range("B53").Value = 2672
$ws.Range("C53").Value = 29
$ws.Range("D53").Value = 1593
$ws.Range("E53").Value = 1033
$ws.Range("G53").Value = 3
$ws.Range("H53").Value = 46

# --- Rows 69/70: Kazajistan overtakes Uzbekistan, swap rows and update Kazajistan's data ---
$ws.Range("A69").Value = "Kazajistan"
$ws.Range("B69").Value = 1331
$ws.Range("C69").Value = 36
$ws.Range("D69").Value = 240
$ws.Range("E69").Value = 1075
$ws.Range("F69").Value = 22
$ws.Range("H69").Value = 16

$ws.Range("A70").Value = "Uzbekistan"
$ws.Range("B70").Value = 1302
$ws.Range("C70").Value = 0
$ws.Range("D70").Value = 107
$ws.Range("E70").Value = 1191
$ws.Range("F70").Value = 8
$ws.Range("H70").Value = 4

# --- Row 98: Kirguistan - update stats (no row swap) ---
$ws.Range("B98").Value = 466
$ws.Range("C98").Value = 17
$ws.Range("D98").Value = 91
$ws.Range("E98").Value = 370

# --- Row 124: El Salvador - update stats (no row swap) ---
$ws.Range("D124").Value = 33
$ws.Range("E124").Value = 120

# --- Row 180: Fiyi - update stats (no row swap) ---
$ws.Range("B180").Value = 17
$ws.Range("C180").Value = 1
$ws.Range("E180").Value = 17
